{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Changes applied (per the authoritative XML diff):\n//   1. \"...event listeners though.\" is split so the cursor / last-edit\n//      marker (_GoBack bookmark) now sits between \"event listeners \" and\n//      \"though.\" instead of after \"all).\" further down the document.\n//   2. The grammar-check markers (proofErr gramStart/gramEnd) that\n//      wrapped the \":hover\" text inside the small floating callout are\n//      cleared (Word drops stale proofing marks once that run is\n//      revisited/re-validated).\n//   3. The \"Background co\" / \"lor\" run split is rejoined into a single\n//      \"Background color\" run.\n\nconst doc = context.document;\nconst body = doc.body;\n\n// ---------------------------------------------------------------------\n// 1. Move the \"_GoBack\" bookmark from after \"all).\" to between\n//    \"event listeners \" and \"though.\" (this mirrors Word automatically\n//    tracking the most recent edit location with the _GoBack bookmark).\n// ---------------------------------------------------------------------\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst thoughResults = body.search(\"though.\", { matchCase: false });\nthoughResults.load(\"text\");\nawait context.sync();\n\nif (thoughResults.items.length > 0) {\n  const thoughStart = thoughResults.items[0].getRange(\"Start\");\n  thoughStart.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2. Drop the stray proofErr (gramStart/gramEnd) markers around the\n//    \":hover\" text inside the floating textbox callout. Office.js has\n//    no direct handle into legacy VML textbox content, so the paragraph\n//    that anchors the callout is round-tripped through Ooxml to let the\n//    host re-normalize it (clearing the now-stale proofing marks).\n// ---------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"pseudo-selector\") !== -1) {\n    const hoverParagraph = paragraphs.items[i];\n    const ooxml = hoverParagraph.getOoxml();\n    await context.sync();\n    hoverParagraph.insertOoxml(ooxml.value, \"Replace\");\n    await context.sync();\n    break;\n  }\n}\n\n// ---------------------------------------------------------------------\n// 3. Rejoin the \"Background co\" + \"lor\" runs into a single\n//    \"Background color\" run.\n// ---------------------------------------------------------------------\nconst bgResults = body.search(\"Background color\", { matchCase: true });\nawait context.sync();\n\nif (bgResults.items.length > 0) {\n  bgResults.items[0].insertText(\"Background color\", \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument / $d / $doc resolve to the open document.\n#\n# Changes applied (per the authoritative XML diff):\n#   1. \"...event listeners though.\" is split so the cursor / last-edit\n#      marker (_GoBack bookmark) now sits between \"event listeners \" and\n#      \"though.\" instead of after \"all).\" further down the document.\n#   2. The grammar-check markers (proofErr gramStart/gramEnd) that\n#      wrapped the \":hover\" text inside the small floating callout are\n#      cleared (Word drops stale proofing marks once that run is\n#      revisited/re-validated).\n#   3. The \"Background co\" / \"lor\" run split is rejoined into a single\n#      \"Background color\" run.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1. Move the \"_GoBack\" bookmark from after \"all).\" to between\n#    \"event listeners \" and \"though.\" (this mirrors Word automatically\n#    tracking the most recent edit location with the _GoBack bookmark).\n# ---------------------------------------------------------------------\ntry {\n    $oldBookmark = $d.Bookmarks(\"_GoBack\")\n    $oldBookmark.Delete()\n} catch {\n    # no existing _GoBack bookmark - nothing to remove\n}\n\n$findRange = $d.Content\nif ($findRange.Find.Execute(\"though.\")) {\n    $newSpot = $findRange.Duplicate\n    $newSpot.Collapse(1)  # wdCollapseStart\n    $d.Bookmarks.Add(\"_GoBack\", $newSpot) | Out-Null\n}\n\n# ---------------------------------------------------------------------\n# 2. Drop the stray proofErr (gramStart/gramEnd) markers around the\n#    \":hover\" text inside the floating textbox callout. The Word object\n#    model has no direct handle into legacy VML textbox content, so the\n#    paragraph that anchors the callout is round-tripped through\n#    WordOpenXML/InsertXML to let the host re-normalize it (clearing the\n#    now-stale proofing marks).\n# ---------------------------------------------------------------------\n$paraCount = $d.Paragraphs.Count\nfor ($i = 1; $i -le $paraCount; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Contains(\"pseudo-selector\")) {\n        $hoverXml = $p.Range.WordOpenXML\n        $p.Range.InsertXML($hoverXml)\n        break\n    }\n}\n\n# ---------------------------------------------------------------------\n# 3. Rejoin the \"Background co\" + \"lor\" runs into a single\n#    \"Background color\" run.\n# ---------------------------------------------------------------------\n$bgRange = $d.Content\n$bgRange.Find.Execute(\"Background color\", $false, $false, $false, $false, $false, $true, 1, $false, \"Background color\", 2) | Out-Null\n"}
